# Sprint 3 burndown update
# - Refresh the "Story Points" (B) and "Guideline" (C) actuals for rows 2-23,
#   including filling in the previously-empty tail rows (12-23).
# - Move the selection to D13 (clears the stale topLeftCell scroll state).
# - Re-home the burndown chart a couple of columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Story Points (col B) / Guideline (col C) values, keyed by row.
$values = @{
    2  = @(19, 19)
    3  = @(19, 19)
    4  = @(19, 19)
    5  = @(19, 19)
    6  = @(19, 19)
    7  = @(19, 19)
    8  = @(16, 19)
    9  = @(16, 19)
    10 = @(16, 19)
    11 = @(16, 19)
    12 = @(16, 16)
    13 = @(13, 16)
    14 = @(13, 16)
    15 = @(13, 16)
    16 = @(13, 13)
    17 = @(9, 13)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(4, 9)
    23 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}

# Move the active selection to D13 (also drops the saved topLeftCell scroll
# position back to the default, matching the target sheetView).
[void]$ws.Range("D13").Select()

# Re-position / re-home the chart: from E5 (col 5, 355600 EMU offset; row 4,
# 0 EMU offset) to column/row 18/30 with the matching offsets, keeping the
# same overall size.
$co = $ws.ChartObjects().Item(1)
$co.Top = 64
$co.Left = 320.1875
$co.Width = 770.6875
$co.Height = 428
